$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2022" column (N) to the table, mirroring the format of
# the existing "2021" column (M), then fill in the reported values. ---

# Row 2 (thin separator row under the title) — blank cell, same style as M2.
$ws.Range("M2").Copy($ws.Range("N2"))

# Row 3 (year header row) — 2022.
$ws.Range("M3").Copy($ws.Range("N3"))
$ws.Range("N3").Value = 2022

# Row 4 — "256 kbit/s to < 2 Mbit/s" subscriptions.
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 6333

# Row 5 — "2 Mbit/s to < 10 Mbit/s" subscriptions.
$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 82675

# Row 6 — ">= 10 Mbit/s" subscriptions.
$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 300853

# Match the selection left behind in the saved file (cell N2 selected).
$ws.Range("N2").Select()
